$d = $word.ActiveDocument

$d.Content.Find.Execute("35×62=2170", $true, $false, $false, $false, $false, $true, 1, $false, "65×38=2470", 2) | Out-Null
$d.Content.Find.Execute("74×65=4810", $true, $false, $false, $false, $false, $true, 1, $false, "64×28=1792", 2) | Out-Null
$d.Content.Find.Execute("87×38=3306", $true, $false, $false, $false, $false, $true, 1, $false, "60×92=5520", 2) | Out-Null
$d.Content.Find.Execute("18×13=234", $true, $false, $false, $false, $false, $true, 1, $false, "75×52=3900", 2) | Out-Null
$d.Content.Find.Execute("86×74=6364", $true, $false, $false, $false, $false, $true, 1, $false, "13×21=273", 2) | Out-Null
$d.Content.Find.Execute("69×42=2898", $true, $false, $false, $false, $false, $true, 1, $false, "93×55=5115", 2) | Out-Null
$d.Content.Find.Execute("51×71=3621", $true, $false, $false, $false, $false, $true, 1, $false, "40×56=2240", 2) | Out-Null
$d.Content.Find.Execute("49×64=3136", $true, $false, $false, $false, $false, $true, 1, $false, "55×87=4785", 2) | Out-Null
$d.Content.Find.Execute("79×65=5135", $true, $false, $false, $false, $false, $true, 1, $false, "39×87=3393", 2) | Out-Null
$d.Content.Find.Execute("31×73=2263", $true, $false, $false, $false, $false, $true, 1, $false, "39×66=2574", 2) | Out-Null
$d.Content.Find.Execute("73×24=1752", $true, $false, $false, $false, $false, $true, 1, $false, "84×14=1176", 2) | Out-Null
$d.Content.Find.Execute("71×91=6461", $true, $false, $false, $false, $false, $true, 1, $false, "24×94=2256", 2) | Out-Null
$d.Content.Find.Execute("75×20=1500", $true, $false, $false, $false, $false, $true, 1, $false, "44×50=2200", 2) | Out-Null
$d.Content.Find.Execute("48×78=3744", $true, $false, $false, $false, $false, $true, 1, $false, "42×86=3612", 2) | Out-Null
$d.Content.Find.Execute("52×22=1144", $true, $false, $false, $false, $false, $true, 1, $false, "34×70=2380", 2) | Out-Null
$d.Content.Find.Execute("12×35=420", $true, $false, $false, $false, $false, $true, 1, $false, "52×97=5044", 2) | Out-Null
$d.Content.Find.Execute("82×96=7872", $true, $false, $false, $false, $false, $true, 1, $false, "38×95=3610", 2) | Out-Null
$d.Content.Find.Execute("45×24=1080", $true, $false, $false, $false, $false, $true, 1, $false, "39×78=3042", 2) | Out-Null
$d.Content.Find.Execute("20×19=380", $true, $false, $false, $false, $false, $true, 1, $false, "33×52=1716", 2) | Out-Null
$d.Content.Find.Execute("14×32=448", $true, $false, $false, $false, $false, $true, 1, $false, "71×13=923", 2) | Out-Null
$d.Content.Find.Execute("61×26=1586", $true, $false, $false, $false, $false, $true, 1, $false, "38×71=2698", 2) | Out-Null
$d.Content.Find.Execute("99×40=3960", $true, $false, $false, $false, $false, $true, 1, $false, "41×54=2214", 2) | Out-Null
$d.Content.Find.Execute("40×34=1360", $true, $false, $false, $false, $false, $true, 1, $false, "74×55=4070", 2) | Out-Null
$d.Content.Find.Execute("83×57=4731", $true, $false, $false, $false, $false, $true, 1, $false, "75×92=6900", 2) | Out-Null
$d.Content.Find.Execute("97×17=1649", $true, $false, $false, $false, $false, $true, 1, $false, "16×33=528", 2) | Out-Null
